$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 289.94116
$ws.Range("I96").Value = 264.53845
$ws.Range("J96").Value = 372.5
$ws.Range("K96").Value = 793.61535
$ws.Range("L96").Value = 1117.5
$ws.Range("M96").Value = 579.38465
$ws.Range("N96").Value = -3863.5

$ws.Range("H106").Value = 6339
$ws.Range("I106").Value = 4578
$ws.Range("K106").Value = 4578
$ws.Range("M106").Value = -3947

$ws.Range("H113").Value = 3139.12
$ws.Range("I113").Value = 2888.3157
$ws.Range("J113").Value = 3933.3333
$ws.Range("K113").Value = 2888.3157
$ws.Range("L113").Value = 3933.3333
$ws.Range("M113").Value = 365.6842999999999
$ws.Range("N113").Value = -10441.3333

$ws.Range("H137").Value = 2327146.8
$ws.Range("I137").Value = 3847441.2
$ws.Range("J137").Value = 1990.6471
$ws.Range("K137").Value = 11542323.6
$ws.Range("L137").Value = 5971.9413
$ws.Range("M137").Value = -11539773.6
$ws.Range("N137").Value = -11071.9413

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1300.6666
$ws.Range("I45").Value = 983.17645
$ws.Range("K45").Value = 983.17645
$ws.Range("M45").Value = -606.17645

$ws.Range("H61").Value = 20876002
$ws.Range("I61").Value = 22750988
$ws.Range("J61").Value = 251153.5
$ws.Range("K61").Value = 22750988
$ws.Range("L61").Value = 251153.5
$ws.Range("M61").Value = -22750776
$ws.Range("N61").Value = -251577.5

$ws.Range("H74").Value = 4667906.5
$ws.Range("I74").Value = 5977248
$ws.Range("J74").Value = 85209.914
$ws.Range("K74").Value = 5977248
$ws.Range("L74").Value = 85209.914
$ws.Range("M74").Value = -5976374
$ws.Range("N74").Value = -86957.914

$ws.Range("H77").Value = 4667906.5
$ws.Range("I77").Value = 5977248
$ws.Range("J77").Value = 85209.914
$ws.Range("K77").Value = 29886240
$ws.Range("L77").Value = 426049.57
$ws.Range("M77").Value = -29881872
$ws.Range("N77").Value = -434785.57

$ws.Range("H132").Value = 105436.8
$ws.Range("I132").Value = 72366.92999999999
$ws.Range("J132").Value = 182599.83
$ws.Range("K132").Value = 217100.79
$ws.Range("L132").Value = 547799.49
$ws.Range("M132").Value = -214570.79
$ws.Range("N132").Value = -552859.49

$ws.Range("H136").Value = 20876002
$ws.Range("I136").Value = 22750988
$ws.Range("J136").Value = 251153.5
$ws.Range("K136").Value = 68252964
$ws.Range("L136").Value = 753460.5
$ws.Range("M136").Value = -68250414
$ws.Range("N136").Value = -758560.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H141").Value = 45151
$ws.Range("J141").Value = 46631.668
$ws.Range("L141").Value = 46631.668
$ws.Range("N141").Value = -56991.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2352.377
$ws.Range("I31").Value = 902.1111
$ws.Range("J31").Value = 6431.25
$ws.Range("K31").Value = 902.1111
$ws.Range("L31").Value = 6431.25
$ws.Range("M31").Value = -607.1111
$ws.Range("N31").Value = -7021.25

$ws.Range("H34").Value = 2352.377
$ws.Range("I34").Value = 902.1111
$ws.Range("J34").Value = 6431.25
$ws.Range("K34").Value = 902.1111
$ws.Range("L34").Value = 6431.25
$ws.Range("M34").Value = -700.1111
$ws.Range("N34").Value = -6835.25

$ws.Range("H64").Value = 23101.834
$ws.Range("J64").Value = 23101.834
$ws.Range("L64").Value = 23101.834
$ws.Range("N64").Value = -23597.834

$ws.Range("H67").Value = 23101.834
$ws.Range("J67").Value = 23101.834
$ws.Range("L67").Value = 23101.834
$ws.Range("N67").Value = -24817.834

$ws.Range("H68").Value = 30397.5
$ws.Range("J68").Value = 32597.143
$ws.Range("L68").Value = 32597.143
$ws.Range("N68").Value = -34095.143

$ws.Range("H71").Value = 30397.5
$ws.Range("J71").Value = 32597.143
$ws.Range("L71").Value = 97791.429
$ws.Range("N71").Value = -105279.429

$ws.Range("H132").Value = 64005.97
$ws.Range("I132").Value = 39952.96
$ws.Range("J132").Value = 168235.67
$ws.Range("K132").Value = 119858.88
$ws.Range("L132").Value = 504707.01
$ws.Range("M132").Value = -117328.88
$ws.Range("N132").Value = -509767.01

$ws.Range("H134").Value = 21736.623
$ws.Range("I134").Value = 1296.3024
$ws.Range("J134").Value = 109630
$ws.Range("K134").Value = 3888.9072
$ws.Range("L134").Value = 328890
$ws.Range("M134").Value = -1353.9072
$ws.Range("N134").Value = -333960

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 1120.7142
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 1120.7142
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3362.1426
$ws.Range("M58").Value = ""
$ws.Range("N58").Value = -3618.1426

$ws.Range("H113").Value = 1056.2858
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 1149
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 3447
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -7787

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 32157.4
$ws.Range("J95").Value = 32157.4
$ws.Range("L95").Value = 32157.4
$ws.Range("N95").Value = -37649.4

$ws.Range("H126").Value = 1550
$ws.Range("I126").Value = 1308.3334
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 3925.0002
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -1455.0002
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 39637.58
$ws.Range("I132").Value = 22364.979
$ws.Range("J132").Value = 202000
$ws.Range("K132").Value = 67094.93700000001
$ws.Range("L132").Value = 606000
$ws.Range("M132").Value = -64564.93700000001
$ws.Range("N132").Value = -611060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2189.1428
$ws.Range("I7").Value = 2134.9285
$ws.Range("J7").Value = 2297.5715
$ws.Range("K7").Value = 2134.9285
$ws.Range("L7").Value = 2297.5715
$ws.Range("M7").Value = -2022.9285
$ws.Range("N7").Value = -2521.5715

$ws.Range("H93").Value = 1019.7692
$ws.Range("I93").Value = 995.3333
$ws.Range("J93").Value = 1074.75
$ws.Range("K93").Value = 995.3333
$ws.Range("L93").Value = 1074.75
$ws.Range("M93").Value = 252.6667
$ws.Range("N93").Value = -3570.75

$ws.Range("H126").Value = 2189.1428
$ws.Range("I126").Value = 2134.9285
$ws.Range("J126").Value = 2297.5715
$ws.Range("K126").Value = 6404.7855
$ws.Range("L126").Value = 6892.7145
$ws.Range("M126").Value = -3934.7855
$ws.Range("N126").Value = -11832.7145

$ws.Range("H132").Value = 32035.559
$ws.Range("I132").Value = 24540.512
$ws.Range("J132").Value = 46699.78
$ws.Range("K132").Value = 73621.53599999999
$ws.Range("L132").Value = 140099.34
$ws.Range("M132").Value = -71091.53599999999
$ws.Range("N132").Value = -145159.34

$ws.Range("H136").Value = 70390.484
$ws.Range("I136").Value = 58953.21
$ws.Range("J136").Value = 88499.5
$ws.Range("K136").Value = 176859.63
$ws.Range("L136").Value = 265498.5
$ws.Range("M136").Value = -174309.63
$ws.Range("N136").Value = -270598.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 40550
$ws.Range("J92").Value = 40550
$ws.Range("L92").Value = 40550
$ws.Range("N92").Value = -45542

$ws.Range("H100").Value = 72743
$ws.Range("I100").Value = 125799.25
$ws.Range("J100").Value = 51520.5
$ws.Range("K100").Value = 251598.5
$ws.Range("L100").Value = 103041
$ws.Range("M100").Value = -251057.5
$ws.Range("N100").Value = -104123

$ws.Range("H105").Value = 20559
$ws.Range("I105").Value = 20559
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 20559
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -17065
$ws.Range("N105").Value = ""

$ws.Range("H122").Value = 3011.5217
$ws.Range("I122").Value = 2697.2144
$ws.Range("J122").Value = 3500.4443
$ws.Range("K122").Value = 8091.6432
$ws.Range("L122").Value = 10501.3329
$ws.Range("M122").Value = -5641.6432
$ws.Range("N122").Value = -15401.3329

$ws.Range("H132").Value = 39367.46
$ws.Range("I132").Value = 33397.605
$ws.Range("J132").Value = 58533.844
$ws.Range("K132").Value = 100192.815
$ws.Range("L132").Value = 175601.532
$ws.Range("M132").Value = -97662.815
$ws.Range("N132").Value = -180661.532

$ws.Range("H136").Value = 36699.055
$ws.Range("I136").Value = 22187.326
$ws.Range("J136").Value = 103453
$ws.Range("K136").Value = 66561.978
$ws.Range("L136").Value = 310359
$ws.Range("M136").Value = -64011.978
$ws.Range("N136").Value = -315459

$ws.Range("H138").Value = 48000
$ws.Range("J138").Value = 48000
$ws.Range("L138").Value = 48000
$ws.Range("N138").Value = -58280
